# adjusted flask api endpoint schema for /stats
#
# Sheet "api_doc": GET /api/health (row 2) and GET /api/stats (row 3) rows
# get their documented error / success JSON payloads updated to match the
# server's actual schema (Python-style True/False booleans, "ok" key
# instead of "status", and a slimmed-down general_error payload without
# the echoed "params" block).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("api_doc")
$ws.Activate()

$ellipsis = [char]0x2026
$lsq = [char]0x201C
$rsq = [char]0x201D

# --- Row 2 (GET /api/health), error column F2 ---
# "false" -> "False" (python boolean casing), rest unchanged.
$cellF2 = $ws.Range("F2")
$cellF2.Value2 = "500: { ""ok"": False, ""error"": ""internal_error"", " + $lsq + "now" + $rsq + ": " + $ellipsis + ", " + $lsq + "response_time_ms" + $rsq + ": 120}"
$cellF2.Characters(1, 3).Font.Bold = $true

# --- Row 3 (GET /api/stats) ---

# D3 (request column) text itself is unchanged.
$ws.Range("D3").Value2 = "Query params only (both mandatory)"

# E3 (success column): "status": "success" -> "ok": True, rest unchanged.
$ws.Range("E3").Value2 = "{""ok"": True, ""num_records"": 120, ""db_connected"": true, ""params"": {""start_date"": " + $ellipsis + ", ""end_date"": " + $ellipsis + "}, ""response_time_ms"": 120, ""now"": " + $ellipsis + ", ""data"": {""day"":{""ave"":{""Fri"":396.25,""Mon"":423.0,""Sat"":360.25,""Sun"":473.25,""Thu"":355.25,""Tue"":327.5,""Wed"":357.25},""std"":{""Fri"":88.205,""Mon"":44.728,""Sat"":131.988,""Sun"":120.477,""Thu"":129.113,""Tue"":147.789,""Wed"":190.902}},""week"":{""ave"":2692.75,""std"":261.545}}"

# F3 (error column): drop echoed "params" block, "status": "failure" -> "ok": True.
$cellF3 = $ws.Range("F3")
$cellF3.Value2 = "400 (general_error): {""ok"": True, ""error"": <python error str>, ""response_time_ms"": " + $ellipsis + ", ""now"": " + $ellipsis + "}"
$cellF3.Characters(1, 3).Font.Bold = $true

# Author's saved cursor position ends on F3.
$ws.Range("F3").Select()
